$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Apply formatting (format-only paste) ---
$ws.Range("A8").Copy() | Out-Null
$ws.Range("A80:N80").PasteSpecial(-4122) | Out-Null
$ws.Range("A9").Copy() | Out-Null
$ws.Range("A81:B81").PasteSpecial(-4122) | Out-Null
$ws.Range("C12").Copy() | Out-Null
$ws.Range("C81").PasteSpecial(-4122) | Out-Null
$ws.Range("D9").Copy() | Out-Null
$ws.Range("D81:E81").PasteSpecial(-4122) | Out-Null
$ws.Range("F9").Copy() | Out-Null
$ws.Range("F81:N81").PasteSpecial(-4122) | Out-Null
$ws.Range("A9").Copy() | Out-Null
$ws.Range("A82:B82").PasteSpecial(-4122) | Out-Null
$ws.Range("C12").Copy() | Out-Null
$ws.Range("C82").PasteSpecial(-4122) | Out-Null
$ws.Range("D9").Copy() | Out-Null
$ws.Range("D82:E82").PasteSpecial(-4122) | Out-Null
$ws.Range("F9").Copy() | Out-Null
$ws.Range("F82:N82").PasteSpecial(-4122) | Out-Null
$ws.Range("A9").Copy() | Out-Null
$ws.Range("A83:B83").PasteSpecial(-4122) | Out-Null
$ws.Range("C9").Copy() | Out-Null
$ws.Range("C83").PasteSpecial(-4122) | Out-Null
$ws.Range("D9").Copy() | Out-Null
$ws.Range("D83:E83").PasteSpecial(-4122) | Out-Null
$ws.Range("F17").Copy() | Out-Null
$ws.Range("F83:N83").PasteSpecial(-4122) | Out-Null
$ws.Range("A9").Copy() | Out-Null
$ws.Range("A84:B84").PasteSpecial(-4122) | Out-Null
$ws.Range("C9").Copy() | Out-Null
$ws.Range("C84").PasteSpecial(-4122) | Out-Null
$ws.Range("D17").Copy() | Out-Null
$ws.Range("D84").PasteSpecial(-4122) | Out-Null
$ws.Range("D9").Copy() | Out-Null
$ws.Range("E84").PasteSpecial(-4122) | Out-Null
$ws.Range("F17").Copy() | Out-Null
$ws.Range("F84:N84").PasteSpecial(-4122) | Out-Null
$ws.Range("A9").Copy() | Out-Null
$ws.Range("A85:B85").PasteSpecial(-4122) | Out-Null
$ws.Range("C9").Copy() | Out-Null
$ws.Range("C85").PasteSpecial(-4122) | Out-Null
$ws.Range("D9").Copy() | Out-Null
$ws.Range("D85").PasteSpecial(-4122) | Out-Null
$ws.Range("D17").Copy() | Out-Null
$ws.Range("E85").PasteSpecial(-4122) | Out-Null
$ws.Range("F9").Copy() | Out-Null
$ws.Range("F85:N85").PasteSpecial(-4122) | Out-Null
$ws.Range("A9").Copy() | Out-Null
$ws.Range("A86:B86").PasteSpecial(-4122) | Out-Null
$ws.Range("C9").Copy() | Out-Null
$ws.Range("C86").PasteSpecial(-4122) | Out-Null
$ws.Range("D17").Copy() | Out-Null
$ws.Range("D86").PasteSpecial(-4122) | Out-Null
$ws.Range("D9").Copy() | Out-Null
$ws.Range("E86").PasteSpecial(-4122) | Out-Null
$ws.Range("F9").Copy() | Out-Null
$ws.Range("F86:N86").PasteSpecial(-4122) | Out-Null
$ws.Range("A15").Copy() | Out-Null
$ws.Range("A87:N87").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Set values ---
$ws.Range("A80").Value2 = "Day"
$ws.Range("B80").Value2 = "Date"
$ws.Range("C80").Value2 = "Ticket"
$ws.Range("D80").Value2 = "Start Time "
$ws.Range("E80").Value2 = "End Time "
$ws.Range("F80").Value2 = "Work Log"
$ws.Range("A81").Value2 = "Friday"
$ws.Range("B81").Value2 = 44995.0
$ws.Range("C81").Value2 = "WBX-4285"
$ws.Range("D81").Value2 = 0.395833333333333
$ws.Range("E81").Value2 = 0.402777777777778
$ws.Range("F81").Value2 = "DSM"
$ws.Range("C82").Value2 = "WBX-4298"
$ws.Range("D82").Value2 = 0.402777777777778
$ws.Range("E82").Value2 = 0.4375
$ws.Range("F82").Value2 = "applied translation on announcement module "
$ws.Range("C83").Value2 = "LMDI-41"
$ws.Range("D83").Value2 = 0.4375
$ws.Range("E83").Value2 = 0.0416666666666667
$ws.Range("F83").Value2 = "R&D for logged user open new tab first open home page  sometimes then dashboard page"
$ws.Range("C84").Value2 = "WBX-4203"
$ws.Range("D84").Value2 = 0.0833333333333333
$ws.Range("E84").Value2 = 0.1875
$ws.Range("F84").Value2 = "tried different approaches N see the changes "
$ws.Range("C85").Value2 = "WBX-4206"
$ws.Range("D85").Value2 = 0.1875
$ws.Range("E85").Value2 = 0.25
$ws.Range("F85").Value2 = "check the code implementation right now is show login page along with dashboard "
$ws.Range("C86").Value2 = "WBX-4284"
$ws.Range("D86").Value2 = 0.25
$ws.Range("E86").Value2 = 0.270833333333333
$ws.Range("F86").Value2 = "ashwani's code working and code implementation is reviewed"

# --- Merge cells ---
$ws.Range("F80:N80").Merge() | Out-Null
$ws.Range("F81:N81").Merge() | Out-Null
$ws.Range("F82:N82").Merge() | Out-Null
$ws.Range("F83:N83").Merge() | Out-Null
$ws.Range("F84:N84").Merge() | Out-Null
$ws.Range("F85:N85").Merge() | Out-Null
$ws.Range("F86:N86").Merge() | Out-Null
$ws.Range("F87:N87").Merge() | Out-Null
$ws.Range("A81:A86").Merge() | Out-Null
$ws.Range("B81:B86").Merge() | Out-Null

Write-Output "DONE"
